$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F6").Value = 1.17

$ws.Range("F9").Value = 1.14
$ws.Range("G9").Value = 1000
$ws.Range("J9").Value = 1.04

$ws.Range("S12").Value = 5.2
